$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume/1h change (E) columns
# D-column values are forced to Text format ("@") before assignment so that
# numeric-looking strings (e.g. "0.9998", "0.06570") keep their exact textual
# representation (trailing zeros, fixed decimal places) instead of being
# auto-converted to a Number by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.431.77"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.14"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.86"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5162"
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3992"
$ws.Range("E8").Value = "  +3.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07903"
$ws.Range("E9").Value = "  -4.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.117"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.92"
$ws.Range("E11").Value = "  -2.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.387"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.46"
$ws.Range("E14").Value = "  -3.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.342"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.807.44"
$ws.Range("E16").Value = "  -1.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.82"
$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001086"
$ws.Range("E18").Value = "  -3.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06570"
$ws.Range("E19").Value = "  -0.94%  "

$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.023"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.465.81"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  -2.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.229"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.26"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.59"
$ws.Range("E27").Value = "  -2.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.020.65"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.402"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.41"
$ws.Range("E30").Value = "  +2.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1095"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.075"
$ws.Range("E32").Value = "  -2.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.666"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.587"
$ws.Range("E34").Value = "  -2.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07250"
$ws.Range("E35").Value = "  -4.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.190"
$ws.Range("E36").Value = "  +4.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02345"
$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2180"
$ws.Range("E38").Value = "  -2.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.67"
$ws.Range("E39").Value = "  -3.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.062"
$ws.Range("E40").Value = "  -3.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6211"
$ws.Range("E41").Value = "  -3.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.158"
$ws.Range("E43").Value = "  -2.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.28"
$ws.Range("E44").Value = "  -2.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6011"
$ws.Range("E45").Value = "  -3.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.312"
$ws.Range("E46").Value = "  -6.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.734"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.90"
$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("E49").Value = "  +1.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.934"
$ws.Range("E50").Value = "  -3.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06844"
$ws.Range("E51").Value = "  -1.88%  "
